$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("ECs", "Vtn", "Itgb8", "ECs", 3, 1, 6.201571333333333, 18.604714, 0.05221490529364391, 0.07406232529850043, 1, 0.3333333333333333, 0.178057, 0.534171, 0.01072706596418309, 0.01486309420510698, 1.104233186899333, 9.938098682093997, 0.0005601127333984909, 0.00110079531796089),
    @("ECs", "Vtn", "Itgb8", "FAPs", 3, 1, 6.201571333333333, 18.604714, 0.05221490529364391, 0.07406232529850043, 3, 1, 2.509764, 7.529292, 0.1512010422647355, 0.2094995353430986, 15.564480475832, 140.080324282488, 0.007894948102153413, 0.01551602273646526),
    @("ECs", "Vtn", "Itgb8", "M2", 3, 1, 6.201571333333333, 18.604714, 0.05221490529364391, 0.07406232529850043, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.003246750280137314, 0.004498597793214309, 0.3342171494864444, 3.007954345377999, 0.0001695287583894817, 0.0003331766131481543),
    @("ECs", "Vtn", "Itgb8", "sCs", 3, 1, 6.201571333333333, 18.604714, 0.05221490529364391, 0.07406232529850043, 2, 1, 13.8571405, 27.714281, 0.8348251414909441, 0.7711387726585802, 85.93604528677233, 515.6162717206339, 0.04359031569970252, 0.05711233063092614),
    @("FAPs", "Vtn", "Itgb8", "ECs", 3, 1, 7.461641333333333, 22.384924, 0.06282422221945559, 0.0891107233935555, 1, 0.3333333333333333, 0.178057, 0.534171, 0.01072706596418309, 0.01486309420510698, 1.328597470889333, 11.957377238004, 0.000673919575896597, 0.001324461076483645),
    @("FAPs", "Vtn", "Itgb8", "FAPs", 3, 1, 7.461641333333333, 22.384924, 0.06282422221945559, 0.0891107233935555, 3, 1, 2.509764, 7.529292, 0.1512010422647355, 0.2094995353430986, 18.726958799312, 168.542629193808, 0.00949908787905304, 0.01866865514503727),
    @("FAPs", "Vtn", "Itgb8", "M2", 3, 1, 7.461641333333333, 22.384924, 0.06282422221945559, 0.0891107233935555, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.003246750280137314, 0.004498597793214309, 0.4021252619497777, 3.619127357547999, 0.0002039745610904263, 0.0004008733036099795),
    @("FAPs", "Vtn", "Itgb8", "sCs", 3, 1, 7.461641333333333, 22.384924, 0.06282422221945559, 0.0891107233935555, 2, 1, 13.8571405, 27.714281, 0.8348251414909441, 0.7711387726585802, 103.3970123166073, 620.382073899644, 0.05244724020341553, 0.06871673386842461),
    @("sCs", "Vtn", "Itgb8", "ECs", 2, 1, 105.106922, 210.213844, 0.8849608724869005, 0.836826951307944, 1, 0.3333333333333333, 0.178057, 0.534171, 0.01072706596418309, 0.01486309420510698, 18.715023210554, 112.290139263324, 0.009493033654888, 0.01243783781066244),
    @("sCs", "Vtn", "Itgb8", "FAPs", 2, 1, 105.106922, 210.213844, 0.8849608724869005, 0.836826951307944, 3, 1, 2.509764, 7.529292, 0.1512010422647355, 0.2094995353430986, 263.793568986408, 1582.761413918448, 0.133807006283529, 0.1753148574615961),
    @("sCs", "Vtn", "Itgb8", "M2", 2, 1, 105.106922, 210.213844, 0.8849608724869005, 0.836826951307944, 1, 0.3333333333333333, 0.05389233333333333, 0.161677, 0.003246750280137314, 0.004498597793214309, 5.664457276064666, 33.98674365638799, 0.002873246960657406, 0.003764547876456175),
    @("sCs", "Vtn", "Itgb8", "sCs", 2, 1, 105.106922, 210.213844, 0.8849608724869005, 0.836826951307944, 2, 1, 13.8571405, 27.714281, 0.8348251414909441, 0.7711387726585802, 1456.481385676541, 5825.925542706163, 0.7387875855878261, 0.6453097081592294)
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

Write-Output "Updated $($data.Count) rows starting at row $startRow"
